# Automatic update of files.
# Bump the "Förändrad" (Changed) date column (C) from 2025-03-01 (45717)
# to 2025-03-02 (45718) for every data row in the sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Avverkningsanmälningar")

$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 3)
    if ($cell.Value2 -eq 45717) {
        $cell.Value2 = 45718
    }
}
